$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the "Meta description" paragraph that currently sits right after
#    the title (Heading1) paragraph.
# ---------------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
if ($metaPara.Range.Text.StartsWith("Meta description")) {
    $metaPara.Range.Delete()
}

# ---------------------------------------------------------------------------
# 2. At the end of the document, turn the last paragraph (the italic
#    "Please create a cartoon style image..." image-prompt paragraph) into
#    two paragraphs:
#      - a new bold paragraph with the page title text
#      - the existing paragraph, now italic, with the meta-description text
# ---------------------------------------------------------------------------
$wordMlNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)

$titleText = "Play Dead or Alive 2 for free - an immersive Western-themed slot game"
$descText  = "Read our review of Dead or Alive 2, a Western-themed slot game with exciting free spins mode, and play for free today on your mobile device."

$replacementXml = '<w:p ' + $wordMlNs + '><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>' + $titleText + '</w:t></w:r></w:p>' + `
                   '<w:p ' + $wordMlNs + '><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>' + $descText + '</w:t></w:r></w:p>'

$lastPara.Range.InsertXML($replacementXml)
